$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date number format from the last existing data row (B230) down the new
# column-B cells so the new date cells share the same style index as the rest of the
# column (reuses existing cellXfs entry instead of registering a new custom format).
$ws.Range("B230").Copy() | Out-Null
$ws.Range("B231:B237").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Append the new workout rows (231-237) pulled from the latest Strava export.
# Row 231
$ws.Cells.Item(231,1).Value = "Eric"
$ws.Cells.Item(231,2).Value = 45491
$ws.Cells.Item(231,3).Value = "Run"
$ws.Cells.Item(231,4).Value = 39
$ws.Cells.Item(231,5).Value = 4
$ws.Cells.Item(231,6).Value = 325
$ws.Cells.Item(231,7).Value = 0
$ws.Cells.Item(231,8).Value = 0
$ws.Cells.Item(231,9).Value = 3
$ws.Cells.Item(231,10).Value = 11
$ws.Cells.Item(231,11).Value = 22
$ws.Cells.Item(231,12).Value = "Wily Hyena"
$ws.Cells.Item(231,13).Value = 6

# Row 232
$ws.Cells.Item(232,1).Value = "Steven"
$ws.Cells.Item(232,2).Value = 45491
$ws.Cells.Item(232,3).Value = "Walk"
$ws.Cells.Item(232,4).Value = 21
$ws.Cells.Item(232,5).Value = 1.03
$ws.Cells.Item(232,6).Value = 23
$ws.Cells.Item(232,7).Value = 21
$ws.Cells.Item(232,8).Value = 0
$ws.Cells.Item(232,9).Value = 0
$ws.Cells.Item(232,10).Value = 0
$ws.Cells.Item(232,11).Value = 0
$ws.Cells.Item(232,12).Value = "Brave Leopard"
$ws.Cells.Item(232,13).Value = 6

# Row 233
$ws.Cells.Item(233,1).Value = "Steven"
$ws.Cells.Item(233,2).Value = 45491
$ws.Cells.Item(233,3).Value = "Walk"
$ws.Cells.Item(233,4).Value = 32
$ws.Cells.Item(233,5).Value = 1.63
$ws.Cells.Item(233,6).Value = 102
$ws.Cells.Item(233,7).Value = 32
$ws.Cells.Item(233,8).Value = 0
$ws.Cells.Item(233,9).Value = 0
$ws.Cells.Item(233,10).Value = 0
$ws.Cells.Item(233,11).Value = 0
$ws.Cells.Item(233,12).Value = "Brave Leopard"
$ws.Cells.Item(233,13).Value = 6

# Row 234
$ws.Cells.Item(234,1).Value = "Phil"
$ws.Cells.Item(234,2).Value = 45491
$ws.Cells.Item(234,3).Value = "Run"
$ws.Cells.Item(234,4).Value = 15
$ws.Cells.Item(234,5).Value = 1.59
$ws.Cells.Item(234,6).Value = 118
$ws.Cells.Item(234,7).Value = 1
$ws.Cells.Item(234,8).Value = 4
$ws.Cells.Item(234,9).Value = 8
$ws.Cells.Item(234,10).Value = 1
$ws.Cells.Item(234,11).Value = 0
$ws.Cells.Item(234,12).Value = "Sauntering Hippo"
$ws.Cells.Item(234,13).Value = 6

# Row 235
$ws.Cells.Item(235,1).Value = "Phil"
$ws.Cells.Item(235,2).Value = 45491
$ws.Cells.Item(235,3).Value = "Workout"
$ws.Cells.Item(235,4).Value = 53
$ws.Cells.Item(235,5).Value = 0
$ws.Cells.Item(235,6).Value = 0
$ws.Cells.Item(235,7).Value = 34
$ws.Cells.Item(235,8).Value = 18
$ws.Cells.Item(235,9).Value = 1
$ws.Cells.Item(235,10).Value = 0
$ws.Cells.Item(235,11).Value = 0
$ws.Cells.Item(235,12).Value = "Sauntering Hippo"
$ws.Cells.Item(235,13).Value = 6

# Row 236
$ws.Cells.Item(236,1).Value = "Phil"
$ws.Cells.Item(236,2).Value = 45491
$ws.Cells.Item(236,3).Value = "Run"
$ws.Cells.Item(236,4).Value = 16
$ws.Cells.Item(236,5).Value = 1.63
$ws.Cells.Item(236,6).Value = 117
$ws.Cells.Item(236,7).Value = 0
$ws.Cells.Item(236,8).Value = 2
$ws.Cells.Item(236,9).Value = 5
$ws.Cells.Item(236,10).Value = 5
$ws.Cells.Item(236,11).Value = 0
$ws.Cells.Item(236,12).Value = "Sauntering Hippo"
$ws.Cells.Item(236,13).Value = 6

# Row 237
$ws.Cells.Item(237,1).Value = "Steven"
$ws.Cells.Item(237,2).Value = 45492
$ws.Cells.Item(237,3).Value = "Walk"
$ws.Cells.Item(237,4).Value = 31
$ws.Cells.Item(237,5).Value = 1.69
$ws.Cells.Item(237,6).Value = 108
$ws.Cells.Item(237,7).Value = 31
$ws.Cells.Item(237,8).Value = 0
$ws.Cells.Item(237,9).Value = 0
$ws.Cells.Item(237,10).Value = 0
$ws.Cells.Item(237,11).Value = 0
$ws.Cells.Item(237,12).Value = "Brave Leopard"
$ws.Cells.Item(237,13).Value = 6

# Move the selection to the first empty row below the newly appended data,
# matching where Excel leaves the cursor after the paste/entry session.
$ws.Range("A238").Select() | Out-Null
